$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2747.5386
$ws.Range("I41").Value = 2997.5
$ws.Range("K41").Value = 2997.5
$ws.Range("M41").Value = -2557.5

$ws.Range("H51").Value = 7002.0605
$ws.Range("J51").Value = 8572.286
$ws.Range("L51").Value = 8572.286
$ws.Range("N51").Value = -9540.286

$ws.Range("H80").Value = 832.1539
$ws.Range("I80").Value = 546.7778
$ws.Range("J80").Value = 1474.25
$ws.Range("K80").Value = 1640.3334
$ws.Range("L80").Value = 4422.75
$ws.Range("M80").Value = -642.3334
$ws.Range("N80").Value = -6418.75

$ws.Range("H83").Value = 832.1539
$ws.Range("I83").Value = 546.7778
$ws.Range("J83").Value = 1474.25
$ws.Range("K83").Value = 4921.000199999999
$ws.Range("L83").Value = 13268.25
$ws.Range("M83").Value = 70.9998000000005
$ws.Range("N83").Value = -23252.25

$ws.Range("H132").Value = 7932.6
$ws.Range("J132").Value = 13269.296
$ws.Range("L132").Value = 39807.888
$ws.Range("N132").Value = -44867.888

$ws.Range("H137").Value = 6947741
$ws.Range("I137").Value = 1026.9231
$ws.Range("K137").Value = 3080.7693
$ws.Range("M137").Value = -530.7692999999999

$ws.Range("H138").Value = 3892.0637
$ws.Range("J138").Value = 4216.974
$ws.Range("L138").Value = 12650.922
$ws.Range("N138").Value = -22930.922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1091890
$ws.Range("I2").Value = 1745784.9
$ws.Range("J2").Value = 2065.3333
$ws.Range("K2").Value = 1745784.9
$ws.Range("L2").Value = 2065.3333
$ws.Range("M2").Value = -1745671.9
$ws.Range("N2").Value = -2291.3333

$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -184

$ws.Range("H61").Value = 17063.637
$ws.Range("I61").Value = 31780.4
$ws.Range("J61").Value = 4799.6665
$ws.Range("K61").Value = 31780.4
$ws.Range("L61").Value = 4799.6665
$ws.Range("M61").Value = -31568.4
$ws.Range("N61").Value = -5223.6665

$ws.Range("H110").Value = 758544.4
$ws.Range("I110").Value = 973000.0600000001
$ws.Range("K110").Value = 973000.0600000001
$ws.Range("M110").Value = -970955.0600000001

$ws.Range("H116").Value = 1091890
$ws.Range("I116").Value = 1745784.9
$ws.Range("J116").Value = 2065.3333
$ws.Range("K116").Value = 1745784.9
$ws.Range("L116").Value = 2065.3333
$ws.Range("M116").Value = -1743490.9
$ws.Range("N116").Value = -6653.3333

$ws.Range("H136").Value = 17063.637
$ws.Range("I136").Value = 31780.4
$ws.Range("J136").Value = 4799.6665
$ws.Range("K136").Value = 95341.20000000001
$ws.Range("L136").Value = 14398.9995
$ws.Range("M136").Value = -92791.20000000001
$ws.Range("N136").Value = -19498.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1091890
$ws.Range("I3").Value = 1745784.9
$ws.Range("J3").Value = 2065.3333
$ws.Range("K3").Value = 1745784.9
$ws.Range("L3").Value = 2065.3333
$ws.Range("M3").Value = -1745670.9
$ws.Range("N3").Value = -2293.3333

$ws.Range("H86").Value = 1360.5
$ws.Range("I86").Value = 1306.9333
$ws.Range("J86").Value = 1521.2
$ws.Range("K86").Value = 1306.9333
$ws.Range("L86").Value = 1521.2
$ws.Range("M86").Value = -183.9332999999999
$ws.Range("N86").Value = -3767.2

$ws.Range("H89").Value = 1360.5
$ws.Range("I89").Value = 1306.9333
$ws.Range("J89").Value = 1521.2
$ws.Range("K89").Value = 6534.666499999999
$ws.Range("L89").Value = 7606
$ws.Range("M89").Value = -918.6664999999994
$ws.Range("N89").Value = -18838

$ws.Range("H94").Value = 721712.1
$ws.Range("I94").Value = 806484.25
$ws.Range("K94").Value = 806484.25
$ws.Range("M94").Value = -806033.25

$ws.Range("H105").Value = 3567.111
$ws.Range("I105").Value = 4028.1428
$ws.Range("K105").Value = 4028.1428
$ws.Range("M105").Value = -2281.1428

$ws.Range("H134").Value = 4477.95
$ws.Range("I134").Value = 2005.1538
$ws.Range("J134").Value = 9070.286
$ws.Range("K134").Value = 6015.4614
$ws.Range("L134").Value = 27210.858
$ws.Range("M134").Value = -3480.4614
$ws.Range("N134").Value = -32280.858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1136777.8
$ws.Range("I107").Value = 1653187.5
$ws.Range("K107").Value = 1653187.5
$ws.Range("M107").Value = -1651267.5

$ws.Range("H132").Value = 11123963
$ws.Range("I132").Value = 12355514
$ws.Range("K132").Value = 37066542
$ws.Range("M132").Value = -37064012

$ws.Range("H140").Value = 39979.91
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 528.6875
$ws.Range("J9").Value = 524.25
$ws.Range("L9").Value = 1572.75
$ws.Range("N9").Value = -2020.75

$ws.Range("H121").Value = 429460
$ws.Range("I121").Value = 1329.3334
$ws.Range("J121").Value = 750558
$ws.Range("K121").Value = 3988.0002
$ws.Range("L121").Value = 2251674
$ws.Range("M121").Value = -2678.0002
$ws.Range("N121").Value = -2254294

$ws.Range("H122").Value = 673.625
$ws.Range("J122").Value = 698.5714
$ws.Range("L122").Value = 6287.1426
$ws.Range("N122").Value = -11187.1426

$ws.Range("H124").Value = 9542.857
$ws.Range("I124").Value = 4450
$ws.Range("K124").Value = 13350
$ws.Range("M124").Value = -8440

$ws.Range("H139").Value = 2947.611
$ws.Range("I139").Value = 1737.4
$ws.Range("J139").Value = 8998.666999999999
$ws.Range("K139").Value = 5212.200000000001
$ws.Range("L139").Value = 26996.001
$ws.Range("M139").Value = -72.20000000000073
$ws.Range("N139").Value = -37276.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 78649.5
$ws.Range("J140").Value = 78649.5
$ws.Range("L140").Value = 78649.5
$ws.Range("N140").Value = -89009.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3730.9092
$ws.Range("I61").Value = 3449.7778
$ws.Range("K61").Value = 3449.7778
$ws.Range("M61").Value = -3247.7778

$ws.Range("H68").Value = 877713.25
$ws.Range("I68").Value = 1339122.4
$ws.Range("K68").Value = 1339122.4
$ws.Range("M68").Value = -1338373.4

$ws.Range("H71").Value = 877713.25
$ws.Range("I71").Value = 1339122.4
$ws.Range("K71").Value = 6695612
$ws.Range("M71").Value = -6691868

$ws.Range("H93").Value = 6493.5
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H113").Value = 3730.9092
$ws.Range("I113").Value = 3449.7778
$ws.Range("K113").Value = 3449.7778
$ws.Range("M113").Value = -1279.7778

$ws.Range("H132").Value = 3853.0967
$ws.Range("I132").Value = 3024.9302
$ws.Range("J132").Value = 5727.3687
$ws.Range("K132").Value = 9074.7906
$ws.Range("L132").Value = 17182.1061
$ws.Range("M132").Value = -6544.7906
$ws.Range("N132").Value = -22242.1061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 21630
$ws.Range("J41").Value = 21630
$ws.Range("L41").Value = 21630
$ws.Range("N41").Value = -22410

$ws.Range("H107").Value = 2772.125
$ws.Range("I107").Value = 3598.182
$ws.Range("J107").Value = 954.8
$ws.Range("K107").Value = 10794.546
$ws.Range("L107").Value = 2864.4
$ws.Range("M107").Value = -8874.545999999998
$ws.Range("N107").Value = -6704.4
